# Apply the changes described by the diff:
# - Insert new shared string "cate" used by A12
# - Rename "UNITS" -> "units" (A20)
# - Remove now-unused "gener_materials" string; A24 becomes numeric 0
# - Add new light-grey font color to header cells G1:G3
# - Add new row: A12 = "cate"
# - Move selection to A21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a light grey font color to the gener_id header column (G1:G3)
$ws.Range("G1:G3").Font.Color = 14540253   # RGB(221,221,221) == 0xDDDDDD

# A6 keeps its existing style unchanged (handled automatically since we
# only touch G1:G3 above)

# New cell A12 with value "cate", using the default/general style
$ws.Range("A12").Value = "cate"

# Rename UNITS -> units
$ws.Range("A20").Value = "units"

# A24 previously held the now-removed "gener_materials" string; it
# becomes a plain numeric zero
$ws.Range("A24").Value = 0

# Update the active selection to A21
$ws.Range("A21").Select() | Out-Null
